$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Target values for rows 48-53 (columns A,B,D,E,F,G,H,P,Q,R)
$data = @{
  48 = @{ A=111974029; B=88180;  D="VU"; E=6276; F="Goliatmusseron";     G="Tricholoma matsutake";   H="(S.Ito & S.Imai) Singer";                  P="Aloppmoarna, Jmt";    Q=439335; R=6952297 }
  49 = @{ A=111974126; B=88180;  D="VU"; E=6276; F="Goliatmusseron";     G="Tricholoma matsutake";   H="(S.Ito & S.Imai) Singer";                  P="Aloppmoarna i S, Jmt"; Q=439290; R=6952209 }
  50 = @{ A=111974134; B=90806;  D="NT"; E=4361; F="Orange taggsvamp";   G="Hydnellum aurantiacum";  H="(Batsch:Fr.) P.Karst.";                     P="Aloppmoarna i S, Jmt"; Q=439400; R=6952207 }
  51 = @{ A=111974133; B=90830;  D="NT"; E=2059; F="Skrovlig taggsvamp"; G="Hydnellum scabrosum";    H="(Fr.) E.Larss., K.H.Larss. & Kõljalg";     P="Aloppmoarna i S, Jmt"; Q=439390; R=6952220 }
  52 = @{ A=111974124; B=90814;  D="LC"; E=4364; F="Dropptaggsvamp";     G="Hydnellum ferrugineum";  H="(Fr.:Fr.) P. Karst.";                       P="Aloppmoarna i S, Jmt"; Q=439276; R=6952197 }
  53 = @{ A=111974125; B=90808;  D="NT"; E=4362; F="Blå taggsvamp";      G="Hydnellum caeruleum";    H="(Hornem.) P.Karst.";                        P="Aloppmoarna i S, Jmt"; Q=439279; R=6952207 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("A$row").Value = $vals.A
  $ws.Range("B$row").Value = $vals.B
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("E$row").Value = $vals.E
  $ws.Range("F$row").Value = $vals.F
  $ws.Range("G$row").Value = $vals.G
  $ws.Range("H$row").Value = $vals.H
  $ws.Range("P$row").Value = $vals.P
  $ws.Range("Q$row").Value = $vals.Q
  $ws.Range("R$row").Value = $vals.R
}
